$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.448.09"
Set-TextValue "E2" "  -0.88%  "
Set-TextValue "D3" "1.626.62"
Set-TextValue "E3" "  -0.49%  "
Set-TextValue "E4" "  +0.24%  "
Set-TextValue "D5" "213.04"
Set-TextValue "E5" "  -0.12%  "
Set-TextValue "E6" "  +1.51%  "
Set-TextValue "E7" "  +0.11%  "
Set-TextValue "E8" "  -1.15%  "
Set-TextValue "D9" "0.0622"
Set-TextValue "E9" "  +0.50%  "
Set-TextValue "D10" "18.94"
Set-TextValue "E10" "  -0.67%  "
Set-TextValue "D11" "0.0844"
Set-TextValue "E11" "  +0.93%  "
Set-TextValue "D12" "1.851.94"
Set-TextValue "E12" "  -0.59%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.14"
Set-TextValue "E13" "  +2.11%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.592.19"
Set-TextValue "E14" "  -2.54%  "
Set-TextValue "D15" "0.522"
Set-TextValue "E15" "  -0.34%  "
Set-TextValue "D16" "64.41"
Set-TextValue "E16" "  +2.28%  "
Set-TextValue "D17" "26.434.35"
Set-TextValue "E17" "  -0.77%  "
Set-TextValue "E18" "  +0.37%  "
Set-TextValue "D19" "215.12"
Set-TextValue "E19" "  +2.61%  "
Set-TextValue "E20" "  +0.19%  "
Set-TextValue "E21" "  -0.07%  "
Set-TextValue "E22" "  +2.10%  "
Set-TextValue "D23" "9.30"
Set-TextValue "E23" "  -0.66%  "
Set-TextValue "D24" "2.00"
Set-TextValue "E24" "  +5.66%  "
Set-TextValue "D25" "147.72"
Set-TextValue "E25" "  +0.95%  "
Set-TextValue "D26" "1.01"
Set-TextValue "E26" "  +0.29%  "
Set-TextValue "E27" "  -0.32%  "
Set-TextValue "D28" "6.83"
Set-TextValue "E28" "  +2.42%  "
Set-TextValue "E29" "  +1.09%  "
Set-TextValue "E30" "  -1.60%  "
Set-TextValue "E31" "  -0.95%  "
Set-TextValue "E32" "  +2.56%  "
Set-TextValue "E33" "  -0.53%  "
Set-TextValue "E34" "  -0.74%  "
Set-TextValue "D35" "1.218.98"
Set-TextValue "E35" "  +4.64%  "
Set-TextValue "E36" "  -1.21%  "
Set-TextValue "E37" "  +3.10%  "
Set-TextValue "E38" "  +0.09%  "
Set-TextValue "D39" "0.794"
Set-TextValue "E39" "  -1.56%  "
Set-TextValue "E40" "  +0.39%  "
Set-TextValue "E41" "  -2.96%  "
Set-TextValue "D42" "0.793"
Set-TextValue "E42" "  -0.08%  "
Set-TextValue "E43" "  +0.36%  "
Set-TextValue "D44" "1.762.30"
Set-TextValue "E44" "  -0.57%  "
Set-TextValue "D45" "92.78"
Set-TextValue "E45" "  +0.54%  "
Set-TextValue "D47" "54.78"
Set-TextValue "E47" "  +0.41%  "
Set-TextValue "E48" "  -0.79%  "
Set-TextValue "E49" "  -0.54%  "
Set-TextValue "D50" "7.52"
Set-TextValue "E50" "  -0.25%  "
Set-TextValue "E51" "  -0.69%  "
